$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = '58.547.04'
$ws.Cells.Item(2, 5).Value2 = '  -3.50%  '

# Row 3
$ws.Cells.Item(3, 4).Value2 = '2.724.12'
$ws.Cells.Item(3, 5).Value2 = '  -6.24%  '

# Row 4
$ws.Cells.Item(4, 5).Value2 = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '503.32'
$ws.Cells.Item(5, 5).Value2 = '  -4.73%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = '141.00'
$ws.Cells.Item(6, 5).Value2 = '  -1.38%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = '0.998'
$ws.Cells.Item(7, 5).Value2 = '  -0.32%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = '0.530'
$ws.Cells.Item(8, 5).Value2 = '  -4.65%  '

# Row 9
$ws.Cells.Item(9, 4).Value2 = '2.736.44'
$ws.Cells.Item(9, 5).Value2 = '  -6.01%  '

# Row 10
$ws.Cells.Item(10, 2).Value2 = 'Toncoin'
$ws.Cells.Item(10, 3).Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = '6.06'
$ws.Cells.Item(10, 5).Value2 = '  +2.19%  '

# Row 11
$ws.Cells.Item(11, 2).Value2 = 'Dogecoin'
$ws.Cells.Item(11, 3).Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '0.105'
$ws.Cells.Item(11, 5).Value2 = '  -2.88%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = '0.348'
$ws.Cells.Item(12, 5).Value2 = '  -3.22%  '

# Row 13
$ws.Cells.Item(13, 5).Value2 = '  +0.98%  '

# Row 14
$ws.Cells.Item(14, 4).Value2 = '3.203.63'
$ws.Cells.Item(14, 5).Value2 = '  -5.99%  '

# Row 15
$ws.Cells.Item(15, 4).Value2 = '58.634.10'
$ws.Cells.Item(15, 5).Value2 = '  -3.23%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '21.73'
$ws.Cells.Item(16, 5).Value2 = '  -3.84%  '

# Row 17
$ws.Cells.Item(17, 2).Value2 = 'ShibaInu'
$ws.Cells.Item(17, 3).Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = '0.0000135'
$ws.Cells.Item(17, 5).Value2 = '  -4.50%  '

# Row 18
$ws.Cells.Item(18, 2).Value2 = 'WrappedEther'
$ws.Cells.Item(18, 3).Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value2 = '2.715.57'
$ws.Cells.Item(18, 5).Value2 = '  -6.54%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = '4.77'
$ws.Cells.Item(19, 5).Value2 = '  -5.04%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = '10.97'
$ws.Cells.Item(20, 5).Value2 = '  -5.67%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = '343.52'
$ws.Cells.Item(21, 5).Value2 = '  -5.50%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '6.27'
$ws.Cells.Item(22, 5).Value2 = '  -4.74%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = '0.998'
$ws.Cells.Item(23, 5).Value2 = '  -0.20%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value2 = '5.69'
$ws.Cells.Item(24, 5).Value2 = '  +0.91%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '62.95'
$ws.Cells.Item(25, 5).Value2 = '  -1.12%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = '0.428'
$ws.Cells.Item(26, 5).Value2 = '  -5.38%  '

# Row 27
$ws.Cells.Item(27, 5).Value2 = '  -4.62%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = '0.994'
$ws.Cells.Item(28, 5).Value2 = '  -0.43%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = '7.51'
$ws.Cells.Item(29, 5).Value2 = '  -4.32%  '

# Row 30
$ws.Cells.Item(30, 4).Value2 = '0.0₃0831'
$ws.Cells.Item(30, 5).Value2 = '  -3.43%  '

# Row 31
$ws.Cells.Item(31, 5).Value2 = '  -0.01%  '

# Row 32
$ws.Cells.Item(32, 2).Value2 = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = '1.60'
$ws.Cells.Item(32, 5).Value2 = '  -4.52%  '

# Row 33
$ws.Cells.Item(33, 2).Value2 = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = '19.16'
$ws.Cells.Item(33, 5).Value2 = '  -2.15%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = '151.29'
$ws.Cells.Item(34, 5).Value2 = '  +1.94%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = '5.43'
$ws.Cells.Item(35, 5).Value2 = '  -2.76%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = '4.20'
$ws.Cells.Item(36, 5).Value2 = '  -3.60%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = '0.951'
$ws.Cells.Item(37, 5).Value2 = '  -5.27%  '

# Row 38
$ws.Cells.Item(38, 5).Value2 = '  -6.02%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '35.93'
$ws.Cells.Item(39, 5).Value2 = '  -5.28%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '1.40'
$ws.Cells.Item(40, 5).Value2 = '  -6.88%  '

# Row 41
$ws.Cells.Item(41, 5).Value2 = '  -3.31%  '

# Row 42
$ws.Cells.Item(42, 4).Value2 = '2.186.59'
$ws.Cells.Item(42, 5).Value2 = '  -6.28%  '

# Row 43
$ws.Cells.Item(43, 5).Value2 = '  -2.52%  '

# Row 44
$ws.Cells.Item(44, 5).Value2 = '  -0.12%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = '0.601'
$ws.Cells.Item(45, 5).Value2 = '  -6.62%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '19.01'
$ws.Cells.Item(46, 5).Value2 = '  -8.46%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = '4.81'
$ws.Cells.Item(47, 5).Value2 = '  -4.78%  '

# Row 48
$ws.Cells.Item(48, 5).Value2 = '  +0.20%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = '0.0228'
$ws.Cells.Item(49, 5).Value2 = '  -3.07%  '

# Row 50
$ws.Cells.Item(50, 5).Value2 = '  -5.42%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = '18.10'
$ws.Cells.Item(51, 5).Value2 = '  -1.83%  '
